# Update the "count" column (F) values on three worksheets to reflect the
# newly generated data snapshot (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 217
$ws1.Range("F7").Value = 117
$ws1.Range("F10").Value = 49
$ws1.Range("F11").Value = 7035
$ws1.Range("F12").Value = 253
$ws1.Range("F13").Value = 403
$ws1.Range("F14").Value = 3482
$ws1.Range("F15").Value = 261
$ws1.Range("F16").Value = 464
$ws1.Range("F18").Value = 584
$ws1.Range("F19").Value = 64

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 49

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 217
$ws4.Range("F9").Value = 117
$ws4.Range("F12").Value = 49
$ws4.Range("F14").Value = 7035
$ws4.Range("F15").Value = 49
$ws4.Range("F16").Value = 253
$ws4.Range("F17").Value = 403
$ws4.Range("F18").Value = 3482
$ws4.Range("F19").Value = 261
$ws4.Range("F20").Value = 464
$ws4.Range("F22").Value = 584
$ws4.Range("F23").Value = 64
